# current meter calibrations.xlsx -- "Need to Proof Report" / Mick finishing manual
# Rework the calibration sheet: rename headers, insert a new amp/avg (J) column
# (shifting the old K "avg/amp" formulas one column right into K, clearing the
# now-unused L column), and add a compact "Verifed Amps / CT % Avg / Amp/CT%Avg"
# summary table (N5:P12) that Mick copy/pasted the six data points + final
# scaling factor into.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 header renames
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "Verifed Current"
$ws.Range("D3").Value = "Percentage Scaled From Kinetis"
$ws.Range("I3").Value = "Average"

# New J3 header ("amp/avg") + K3 becomes "avg/amp"; old L3 is cleared since the
# avg/amp column now lives in K.
$ws.Range("J3").Value = "amp/avg"
$ws.Range("K3").Value = "avg/amp"
$ws.Range("L3").ClearContents()

# ---------------------------------------------------------------------------
# Shift the "avg/amp" formulas from L into K (K used to hold "amp/avg"); add
# the new "amp/avg" formulas into J. Clear L afterwards.
# ---------------------------------------------------------------------------
$ws.Range("J6").Formula = "=C6/I6"
$ws.Range("J7").Formula = "=C7/I7"
$ws.Range("J9").Formula = "=C9/I9"
$ws.Range("J10").Formula = "=C10/I10"
$ws.Range("J15").Formula = "=C15/I15"
$ws.Range("J17").Formula = "=C17/I17"

$ws.Range("K6").Formula = "=I6/C6"
$ws.Range("K7").Formula = "=I7/C7"
$ws.Range("K9").Formula = "=I9/C9"
$ws.Range("K10").Formula = "=I10/C10"
$ws.Range("K15").Formula = "=I15/C15"
$ws.Range("K17").Formula = "=I17/C17"

$ws.Range("L6").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("L17").ClearContents()

# The summary "SUM(K6:K17)/6" cell moves from K20 to J20 (same calculation,
# now over the relocated J column).
$ws.Range("K20").ClearContents()
$ws.Range("J20").Formula = "=SUM(J6:J17)/6"

# ---------------------------------------------------------------------------
# New "Verifed Amps / CT % Avg / Amp/CT%Avg" summary block, N5:P12
# (O5/P5 typed before N5, matching the shared-string allocation order)
# ---------------------------------------------------------------------------
$ws.Range("O5").Value = "CT % Avg"
$ws.Range("P5").Value = "Amp/CT%Avg"
$ws.Range("N5").Value = "Verifed Amps"

$ws.Range("N6").Value = 0.2
$ws.Range("O6").Value = 1.3679999999999999
$ws.Range("P6").Value = 0.14619883040935674

$ws.Range("N7").Value = 0.2
$ws.Range("O7").Value = 1.2819999999999998
$ws.Range("P7").Value = 0.15600624024961002

$ws.Range("N8").Value = 0.5
$ws.Range("O8").Value = 3.3379999999999996
$ws.Range("P8").Value = 0.14979029358897544

$ws.Range("N9").Value = 0.5
$ws.Range("O9").Value = 3.3220000000000001
$ws.Range("P9").Value = 0.15051173991571343

$ws.Range("N10").Value = 8.8000000000000007
$ws.Range("O10").Value = 48.808000000000007
$ws.Range("P10").Value = 0.18029831175217176

$ws.Range("N11").Value = 8.8000000000000007
$ws.Range("O11").Value = 49.044000000000004
$ws.Range("P11").Value = 0.17943071527607862

$ws.Range("N12").Value = "Scaling Factor"
$ws.Range("O12").Value = "'="
$ws.Range("P12").Value = 0.1603726885319843

# ---------------------------------------------------------------------------
# Cosmetics: column widths for the touched columns + selection left where
# Mick was last working (the new summary block) + narrower window.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.5546875
$ws.Columns.Item(14).ColumnWidth = 13.5546875
$ws.Columns.Item(15).ColumnWidth = 8.33203125
$ws.Columns.Item(16).ColumnWidth = 12

$ws.Range("N5:P12").Select()
